# Trix.Calc.xlsx edit: rename "index" column to "i" and switch from
# 1-based row index to 0-based row index (shift values down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRIX")

# Rename the header in column A (feeds sharedStrings.xml + the
# "testdata" table's first column definition automatically).
$ws.Range("A1").Value2 = "i"

# Shrink column A now that values are shorter ("i" + 0-based numbers).
$ws.Columns.Item(1).ColumnWidth = 3.14

# Re-number the data rows: old column was 1..502 (row 2..503),
# new column is 0..501 (row 2..503).
$lastRow = 503
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
